# Auto-generated edit applying per-cell numeric updates described by the commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 475.9091
$ws.Range("J12").Value2 = 748.75
$ws.Range("L12").Value2 = 748.75
$ws.Range("N12").Value2 = -1088.75
$ws.Range("H19").Value2 = 1984.8125
$ws.Range("I19").Value2 = 1228.6364
$ws.Range("K19").Value2 = 1228.6364
$ws.Range("M19").Value2 = -1053.6364
$ws.Range("H32").Value2 = 4993.1113
$ws.Range("I32").Value2 = 1848.4286
$ws.Range("J32").Value2 = 15999.5
$ws.Range("K32").Value2 = 1848.4286
$ws.Range("L32").Value2 = 15999.5
$ws.Range("M32").Value2 = -1522.4286
$ws.Range("N32").Value2 = -16651.5
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("K53").Value2 = 0
$ws.Range("M53").ClearContents()
$ws.Range("H103").Value2 = 520.9091
$ws.Range("J103").Value2 = 520.9091
$ws.Range("L103").Value2 = 1562.7273
$ws.Range("N103").Value2 = -2734.7273
$ws.Range("H132").Value2 = 2833.5688
$ws.Range("I132").Value2 = 2551.7637
$ws.Range("J132").Value2 = 8000
$ws.Range("K132").Value2 = 7655.2911
$ws.Range("L132").Value2 = 24000
$ws.Range("M132").Value2 = -5125.2911
$ws.Range("N132").Value2 = -29060
$ws.Range("H137").Value2 = 2031.4324
$ws.Range("J137").Value2 = 2358.2222
$ws.Range("L137").Value2 = 7074.6666
$ws.Range("N137").Value2 = -12174.6666
$ws.Range("H138").Value2 = 4447.9
$ws.Range("J138").Value2 = 5628.4062
$ws.Range("L138").Value2 = 16885.2186
$ws.Range("N138").Value2 = -27165.2186

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 601
$ws.Range("I4").Value2 = 477
$ws.Range("K4").Value2 = 477
$ws.Range("M4").Value2 = -361
$ws.Range("H6").Value2 = 5078200
$ws.Range("I6").Value2 = 150000
$ws.Range("J6").Value2 = 8363666.5
$ws.Range("K6").Value2 = 150000
$ws.Range("L6").Value2 = 8363666.5
$ws.Range("M6").Value2 = -149827
$ws.Range("N6").Value2 = -8364012.5
$ws.Range("H32").Value2 = 10171.172
$ws.Range("I32").Value2 = 4990.909
$ws.Range("K32").Value2 = 4990.909
$ws.Range("M32").Value2 = -4703.909
$ws.Range("H63").Value2 = 3979.8
$ws.Range("I63").Value2 = 3983
$ws.Range("K63").Value2 = 3983
$ws.Range("M63").Value2 = -3297
$ws.Range("H66").Value2 = 3979.8
$ws.Range("I66").Value2 = 3983
$ws.Range("K66").Value2 = 19915
$ws.Range("M66").Value2 = -16483
$ws.Range("H74").Value2 = 2528.4707
$ws.Range("I74").Value2 = 2043.3334
$ws.Range("J74").Value2 = 3692.8
$ws.Range("K74").Value2 = 2043.3334
$ws.Range("L74").Value2 = 3692.8
$ws.Range("M74").Value2 = -1169.3334
$ws.Range("N74").Value2 = -5440.8
$ws.Range("H77").Value2 = 2528.4707
$ws.Range("I77").Value2 = 2043.3334
$ws.Range("J77").Value2 = 3692.8
$ws.Range("K77").Value2 = 10216.667
$ws.Range("L77").Value2 = 18464
$ws.Range("M77").Value2 = -5848.666999999999
$ws.Range("N77").Value2 = -27200
$ws.Range("H97").Value2 = 797.7143
$ws.Range("I97").Value2 = 733.6667
$ws.Range("K97").Value2 = 733.6667
$ws.Range("M97").Value2 = -237.6667
$ws.Range("H111").Value2 = 50000
$ws.Range("J111").Value2 = 50000
$ws.Range("L111").Value2 = 50000
$ws.Range("N111").Value2 = -58180
$ws.Range("H122").Value2 = 19735
$ws.Range("I122").Value2 = 20576.875
$ws.Range("K122").Value2 = 61730.625
$ws.Range("M122").Value2 = -59280.625

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value2 = 11666.667
$ws.Range("J3").Value2 = 15000
$ws.Range("L3").Value2 = 15000
$ws.Range("N3").Value2 = -15226
$ws.Range("H31").Value2 = 5952.304
$ws.Range("I31").Value2 = 4375.385
$ws.Range("J31").Value2 = 8002.3
$ws.Range("K31").Value2 = 4375.385
$ws.Range("L31").Value2 = 8002.3
$ws.Range("M31").Value2 = -4080.385
$ws.Range("N31").Value2 = -8592.299999999999
$ws.Range("H34").Value2 = 5952.304
$ws.Range("I34").Value2 = 4375.385
$ws.Range("J34").Value2 = 8002.3
$ws.Range("K34").Value2 = 4375.385
$ws.Range("L34").Value2 = 8002.3
$ws.Range("M34").Value2 = -4173.385
$ws.Range("N34").Value2 = -8406.299999999999
$ws.Range("H58").Value2 = 3151.5557
$ws.Range("I58").Value2 = 4077
$ws.Range("J58").Value2 = 1994.75
$ws.Range("K58").Value2 = 4077
$ws.Range("L58").Value2 = 1994.75
$ws.Range("M58").Value2 = -3874
$ws.Range("N58").Value2 = -2400.75
$ws.Range("H99").Value2 = 4768.8
$ws.Range("I99").Value2 = 4114.684
$ws.Range("K99").Value2 = 4114.684
$ws.Range("M99").Value2 = -2616.684
$ws.Range("H126").Value2 = 4768.8
$ws.Range("I126").Value2 = 4114.684
$ws.Range("K126").Value2 = 12344.052
$ws.Range("M126").Value2 = -9874.052
$ws.Range("H129").Value2 = 89000
$ws.Range("J129").Value2 = 89000
$ws.Range("L129").Value2 = 89000
$ws.Range("N129").Value2 = -99000
$ws.Range("H132").Value2 = 3250.5
$ws.Range("I132").Value2 = 2960
$ws.Range("K132").Value2 = 8880
$ws.Range("M132").Value2 = -6350
$ws.Range("H136").Value2 = 3151.5557
$ws.Range("I136").Value2 = 4077
$ws.Range("J136").Value2 = 1994.75
$ws.Range("K136").Value2 = 12231
$ws.Range("L136").Value2 = 5984.25
$ws.Range("M136").Value2 = -9681
$ws.Range("N136").Value2 = -11084.25
$ws.Range("H141").Value2 = 529850
$ws.Range("I141").Value2 = 0
$ws.Range("K141").Value2 = 0
$ws.Range("M141").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value2 = 781.3333
$ws.Range("J98").Value2 = 781.3333
$ws.Range("L98").Value2 = 2343.9999
$ws.Range("N98").Value2 = -5339.9999
$ws.Range("H131").Value2 = 2439.353
$ws.Range("J131").Value2 = 2870.6667
$ws.Range("L131").Value2 = 8612.000100000001
$ws.Range("N131").Value2 = -18692.0001
$ws.Range("H132").Value2 = 1446.5
$ws.Range("I132").Value2 = 895
$ws.Range("K132").Value2 = 8055
$ws.Range("M132").Value2 = -5525
$ws.Range("H137").Value2 = 4715.353
$ws.Range("J137").Value2 = 4123.1113
$ws.Range("L137").Value2 = 12369.3339
$ws.Range("N137").Value2 = -22569.3339

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 9900
$ws.Range("I126").Value2 = 8998.429
$ws.Range("J126").Value2 = 12003.667
$ws.Range("K126").Value2 = 26995.287
$ws.Range("L126").Value2 = 36011.001
$ws.Range("M126").Value2 = -24525.287
$ws.Range("N126").Value2 = -40951.001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 1517.625
$ws.Range("I46").Value2 = 422
$ws.Range("J46").Value2 = 1882.8334
$ws.Range("K46").Value2 = 422
$ws.Range("L46").Value2 = 1882.8334
$ws.Range("M46").Value2 = -234
$ws.Range("N46").Value2 = -2258.8334
$ws.Range("H122").Value2 = 7197.25
$ws.Range("I122").Value2 = 6663
$ws.Range("J122").Value2 = 8800
$ws.Range("K122").Value2 = 19989
$ws.Range("L122").Value2 = 26400
$ws.Range("M122").Value2 = -17539
$ws.Range("N122").Value2 = -31300
$ws.Range("H136").Value2 = 4156.5557
$ws.Range("I136").Value2 = 4254
$ws.Range("K136").Value2 = 12762
$ws.Range("M136").Value2 = -10212

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value2 = 9881.333000000001
$ws.Range("J41").Value2 = 9881.333000000001
$ws.Range("L41").Value2 = 9881.333000000001
$ws.Range("N41").Value2 = -10661.333
$ws.Range("H74").Value2 = 4173.7144
$ws.Range("J74").Value2 = 3889.3333
$ws.Range("L74").Value2 = 3889.3333
$ws.Range("N74").Value2 = -5761.3333
$ws.Range("H77").Value2 = 4173.7144
$ws.Range("J77").Value2 = 3889.3333
$ws.Range("L77").Value2 = 11667.9999
$ws.Range("N77").Value2 = -21027.9999
$ws.Range("H96").Value2 = 2200
$ws.Range("I96").Value2 = 0
$ws.Range("J96").Value2 = 2200
$ws.Range("K96").Value2 = 0
$ws.Range("L96").Value2 = 2200
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value2 = -4946
$ws.Range("H108").Value2 = 0
$ws.Range("J108").Value2 = 0
$ws.Range("L108").Value2 = 0
$ws.Range("N108").ClearContents()
$ws.Range("H117").Value2 = 32000
$ws.Range("J117").Value2 = 32000
$ws.Range("L117").Value2 = 32000
$ws.Range("N117").Value2 = -41178
$ws.Range("H126").Value2 = 2703.0625
$ws.Range("I126").Value2 = 2269
$ws.Range("J126").Value2 = 5741.5
$ws.Range("K126").Value2 = 6807
$ws.Range("L126").Value2 = 17224.5
$ws.Range("M126").Value2 = -4337
$ws.Range("N126").Value2 = -22164.5
$ws.Range("H127").Value2 = 0
$ws.Range("I127").Value2 = 0
$ws.Range("K127").Value2 = 0
$ws.Range("M127").ClearContents()
$ws.Range("H132").Value2 = 11471.3
$ws.Range("J132").Value2 = 13307.728
$ws.Range("L132").Value2 = 39923.18399999999
$ws.Range("N132").Value2 = -44983.18399999999
$ws.Range("H136").Value2 = 7672.9165
$ws.Range("I136").Value2 = 8791.777
$ws.Range("J136").Value2 = 4316.3335
$ws.Range("K136").Value2 = 26375.331
$ws.Range("L136").Value2 = 12949.0005
$ws.Range("M136").Value2 = -23825.331
$ws.Range("N136").Value2 = -18049.0005
